$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "CV_Santiago",
    "CV_SantiagoRodriguez",
    "CV_SantiagoRodriguez",
    "CV_SRF",
    "GeneralCV",
    "LACCD_CL",
    "Resume_Santiago",
    "RodriguezFlores_Santiago_Resume",
    "RodriguezSantiago_Resume2021 - Copy",
    "SBCC_Resume",
    "SMC_CL",
    "tax2022"
)

$row = 14
foreach ($value in $values) {
    $ws.Cells.Item($row, 2).Value = $value
    $row = $row + 1
}
